$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New weekly-total column (F) ----------------------------------------
# Week totals: formulas summing 5-row blocks of column D, formatted as
# plain integers (numFmtId 1 / "0"), mirroring the style used elsewhere.
$weekRanges = @(
    @{ Cell = "F5";  Formula = "=SUM(D3:D5)" },
    @{ Cell = "F10"; Formula = "=SUM(D6:D10)" },
    @{ Cell = "F15"; Formula = "=SUM(D11:D15)" },
    @{ Cell = "F20"; Formula = "=SUM(D16:D20)" },
    @{ Cell = "F25"; Formula = "=SUM(D21:D25)" },
    @{ Cell = "F30"; Formula = "=SUM(D26:D30)" },
    @{ Cell = "F35"; Formula = "=SUM(D31:D35)" }
)
foreach ($w in $weekRanges) {
    $r = $ws.Range($w.Cell)
    $r.Formula = $w.Formula
    $r.NumberFormat = "0"
}

# Column F width (new column introduced next to the description column)
$ws.Columns.Item(6).ColumnWidth = 11.52

# --- Fill in the previously-blank days (18 Aug - 21 Aug / rows 18-21) ---
$ws.Range("C18").Value2 = 13

$ws.Range("B19").Value2 = 11
$ws.Range("C19").Value2 = 13
# New shared-string entries are appended in first-use order, so write these
# description cells in the same order the author originally typed them:
# SQLFlite (E19), Hive (E18), Raider.io (E21), content-management (E20).
$ws.Range("E19").Value2 = "Looked into SQLFlite options, minor style changes. "

$ws.Range("E18").Value2 = "Spent a few hours playing with Hive but decided not to use it and reverted the changes (lack of documentation, many deprecations)"

$ws.Range("B20").Value2 = 8
$ws.Range("C20").Value2 = 11

$ws.Range("B21").Value2 = 9
$ws.Range("C21").Value2 = 13
$ws.Range("E21").Value2 = "Planning out Raider.io API integration"

$ws.Range("E20").Value2 = "Planned out content management more // had job interview on this day"

# --- Selection moves to E23 ----------------------------------------------
$ws.Range("E23").Select() | Out-Null
